$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 71,6

$arr[0,0] = "Hubei"
$arr[0,1] = "Mainland China"
$arr[0,2] = "2/6/20 23:23"
$arr[0,3] = 22112
$arr[0,4] = 817
$arr[0,5] = 618
$arr[1,0] = "Guangdong"
$arr[1,1] = "Mainland China"
$arr[1,2] = "2/6/20 12:43"
$arr[1,3] = 970
$arr[1,4] = 69
$arr[1,5] = 0
$arr[2,0] = "Zhejiang"
$arr[2,1] = "Mainland China"
$arr[2,2] = "2/6/20 10:53"
$arr[2,3] = 954
$arr[2,4] = 94
$arr[2,5] = 0
$arr[3,0] = "Henan"
$arr[3,1] = "Mainland China"
$arr[3,2] = "2/6/20 11:23"
$arr[3,3] = 851
$arr[3,4] = 56
$arr[3,5] = 2
$arr[4,0] = "Hunan"
$arr[4,1] = "Mainland China"
$arr[4,2] = "2/6/20 13:13"
$arr[4,3] = 711
$arr[4,4] = 81
$arr[4,5] = 0
$arr[5,0] = "Jiangxi"
$arr[5,1] = "Mainland China"
$arr[5,2] = "2/6/20 2:33"
$arr[5,3] = 600
$arr[5,4] = 37
$arr[5,5] = 0
$arr[6,0] = "Anhui"
$arr[6,1] = "Mainland China"
$arr[6,2] = "2/6/20 13:33"
$arr[6,3] = 591
$arr[6,4] = 34
$arr[6,5] = 0
$arr[7,0] = "Chongqing"
$arr[7,1] = "Mainland China"
$arr[7,2] = "2/6/20 23:33"
$arr[7,3] = 411
$arr[7,4] = 24
$arr[7,5] = 2
$arr[8,0] = "Jiangsu"
$arr[8,1] = "Mainland China"
$arr[8,2] = "2/6/20 8:03"
$arr[8,3] = 373
$arr[8,4] = 34
$arr[8,5] = 0
$arr[9,0] = "Shandong"
$arr[9,1] = "Mainland China"
$arr[9,2] = "2/6/20 7:53"
$arr[9,3] = 347
$arr[9,4] = 27
$arr[9,5] = 0
$arr[10,0] = "Sichuan"
$arr[10,1] = "Mainland China"
$arr[10,2] = "2/7/20 0:53"
$arr[10,3] = 344
$arr[10,4] = 37
$arr[10,5] = 1
$arr[11,0] = "Beijing"
$arr[11,1] = "Mainland China"
$arr[11,2] = "2/6/20 3:23"
$arr[11,3] = 274
$arr[11,4] = 31
$arr[11,5] = 1
$arr[12,0] = "Shanghai"
$arr[12,1] = "Mainland China"
$arr[12,2] = "2/7/20 0:03"
$arr[12,3] = 269
$arr[12,4] = 25
$arr[12,5] = 1
$arr[13,0] = "Heilongjiang"
$arr[13,1] = "Mainland China"
$arr[13,2] = "2/6/20 8:13"
$arr[13,3] = 227
$arr[13,4] = 8
$arr[13,5] = 3
$arr[14,0] = "Fujian"
$arr[14,1] = "Mainland China"
$arr[14,2] = "2/6/20 11:03"
$arr[14,3] = 215
$arr[14,4] = 14
$arr[14,5] = 0
$arr[15,0] = "Shaanxi"
$arr[15,1] = "Mainland China"
$arr[15,2] = "2/6/20 13:53"
$arr[15,3] = 173
$arr[15,4] = 9
$arr[15,5] = 0
$arr[16,0] = "Guangxi"
$arr[16,1] = "Mainland China"
$arr[16,2] = "2/7/20 0:33"
$arr[16,3] = 172
$arr[16,4] = 17
$arr[16,5] = 0
$arr[17,0] = "Hebei"
$arr[17,1] = "Mainland China"
$arr[17,2] = "2/7/20 0:43"
$arr[17,3] = 171
$arr[17,4] = 16
$arr[17,5] = 1
$arr[18,0] = "Yunnan"
$arr[18,1] = "Mainland China"
$arr[18,2] = "2/6/20 9:43"
$arr[18,3] = 133
$arr[18,4] = 7
$arr[18,5] = 0
$arr[19,0] = "Hainan"
$arr[19,1] = "Mainland China"
$arr[19,2] = "2/7/20 0:03"
$arr[19,3] = 106
$arr[19,4] = 8
$arr[19,5] = 2
$arr[20,0] = "Shanxi"
$arr[20,1] = "Mainland China"
$arr[20,2] = "2/6/20 23:03"
$arr[20,3] = 96
$arr[20,4] = 12
$arr[20,5] = 0
$arr[21,0] = "Liaoning"
$arr[21,1] = "Mainland China"
$arr[21,2] = "2/6/20 14:33"
$arr[21,3] = 94
$arr[21,4] = 5
$arr[21,5] = 0
$arr[22,0] = "Tianjin"
$arr[22,1] = "Mainland China"
$arr[22,2] = "2/6/20 23:43"
$arr[22,3] = 79
$arr[22,4] = 2
$arr[22,5] = 1
$arr[23,0] = "Guizhou"
$arr[23,1] = "Mainland China"
$arr[23,2] = "2/6/20 9:53"
$arr[23,3] = 71
$arr[23,4] = 6
$arr[23,5] = 1
$arr[24,0] = "Jilin"
$arr[24,1] = "Mainland China"
$arr[24,2] = "2/7/20 0:43"
$arr[24,3] = 65
$arr[24,4] = 4
$arr[24,5] = 1
$arr[25,0] = "Gansu"
$arr[25,1] = "Mainland China"
$arr[25,2] = "2/5/20 16:23"
$arr[25,3] = 62
$arr[25,4] = 6
$arr[25,5] = 0
$arr[26,0] = "Inner Mongolia"
$arr[26,1] = "Mainland China"
$arr[26,2] = "2/6/20 2:23"
$arr[26,3] = 46
$arr[26,4] = 4
$arr[26,5] = 0
$arr[27,0] = ""
$arr[27,1] = "Japan"
$arr[27,2] = "2/6/20 2:53"
$arr[27,3] = 45
$arr[27,4] = 1
$arr[27,5] = 0
$arr[28,0] = "Ningxia"
$arr[28,1] = "Mainland China"
$arr[28,2] = "2/6/20 2:13"
$arr[28,3] = 40
$arr[28,4] = 1
$arr[28,5] = 0
$arr[29,0] = "Xinjiang"
$arr[29,1] = "Mainland China"
$arr[29,2] = "2/6/20 1:13"
$arr[29,3] = 36
$arr[29,4] = 0
$arr[29,5] = 0
$arr[30,0] = ""
$arr[30,1] = "Singapore"
$arr[30,2] = "2/5/20 16:33"
$arr[30,3] = 28
$arr[30,4] = 0
$arr[30,5] = 0
$arr[31,0] = ""
$arr[31,1] = "Thailand"
$arr[31,2] = "2/4/20 15:33"
$arr[31,3] = 25
$arr[31,4] = 5
$arr[31,5] = 0
$arr[32,0] = "Hong Kong"
$arr[32,1] = "Hong Kong"
$arr[32,2] = "2/6/20 14:43"
$arr[32,3] = 24
$arr[32,4] = 0
$arr[32,5] = 1
$arr[33,0] = ""
$arr[33,1] = "South Korea"
$arr[33,2] = "2/6/20 2:53"
$arr[33,3] = 23
$arr[33,4] = 0
$arr[33,5] = 0
$arr[34,0] = "Qinghai"
$arr[34,1] = "Mainland China"
$arr[34,2] = "2/6/20 2:13"
$arr[34,3] = 18
$arr[34,4] = 3
$arr[34,5] = 0
$arr[35,0] = "Taiwan"
$arr[35,1] = "Taiwan"
$arr[35,2] = "2/6/20 15:03"
$arr[35,3] = 16
$arr[35,4] = 1
$arr[35,5] = 0
$arr[36,0] = ""
$arr[36,1] = "Germany"
$arr[36,2] = "2/3/20 20:53"
$arr[36,3] = 12
$arr[36,4] = 0
$arr[36,5] = 0
$arr[37,0] = ""
$arr[37,1] = "Malaysia"
$arr[37,2] = "2/5/20 15:43"
$arr[37,3] = 12
$arr[37,4] = 0
$arr[37,5] = 0
$arr[38,0] = "Macau"
$arr[38,1] = "Macau"
$arr[38,2] = "2/6/20 14:23"
$arr[38,3] = 10
$arr[38,4] = 1
$arr[38,5] = 0
$arr[39,0] = ""
$arr[39,1] = "Vietnam"
$arr[39,2] = "2/6/20 1:13"
$arr[39,3] = 10
$arr[39,4] = 1
$arr[39,5] = 0
$arr[40,0] = ""
$arr[40,1] = "France"
$arr[40,2] = "2/1/20 1:52"
$arr[40,3] = 6
$arr[40,4] = 0
$arr[40,5] = 0
$arr[41,0] = ""
$arr[41,1] = "United Arab Emirates"
$arr[41,2] = "2/2/20 5:43"
$arr[41,3] = 5
$arr[41,4] = 0
$arr[41,5] = 0
$arr[42,0] = "New South Wales"
$arr[42,1] = "Australia"
$arr[42,2] = "2/6/20 3:13"
$arr[42,3] = 4
$arr[42,4] = 2
$arr[42,5] = 0
$arr[43,0] = "Queensland"
$arr[43,1] = "Australia"
$arr[43,2] = "2/6/20 2:53"
$arr[43,3] = 4
$arr[43,4] = 0
$arr[43,5] = 0
$arr[44,0] = "Victoria"
$arr[44,1] = "Australia"
$arr[44,2] = "2/1/20 18:12"
$arr[44,3] = 4
$arr[44,4] = 0
$arr[44,5] = 0
$arr[45,0] = ""
$arr[45,1] = "India"
$arr[45,2] = "2/3/20 21:43"
$arr[45,3] = 3
$arr[45,4] = 0
$arr[45,5] = 0
$arr[46,0] = "South Australia"
$arr[46,1] = "Australia"
$arr[46,2] = "2/2/20 22:33"
$arr[46,3] = 2
$arr[46,4] = 0
$arr[46,5] = 0
$arr[47,0] = "British Columbia"
$arr[47,1] = "Canada"
$arr[47,2] = "2/5/20 17:33"
$arr[47,3] = 2
$arr[47,4] = 0
$arr[47,5] = 0
$arr[48,0] = "Toronto, ON"
$arr[48,1] = "Canada"
$arr[48,2] = "2/4/20 0:13"
$arr[48,3] = 2
$arr[48,4] = 0
$arr[48,5] = 0
$arr[49,0] = ""
$arr[49,1] = "Italy"
$arr[49,2] = "1/31/20 8:15"
$arr[49,3] = 2
$arr[49,4] = 0
$arr[49,5] = 0
$arr[50,0] = ""
$arr[50,1] = "Philippines"
$arr[50,2] = "2/2/20 3:33"
$arr[50,3] = 2
$arr[50,4] = 0
$arr[50,5] = 1
$arr[51,0] = ""
$arr[51,1] = "Russia"
$arr[51,2] = "1/31/20 16:13"
$arr[51,3] = 2
$arr[51,4] = 0
$arr[51,5] = 0
$arr[52,0] = ""
$arr[52,1] = "UK"
$arr[52,2] = "2/1/20 1:52"
$arr[52,3] = 2
$arr[52,4] = 0
$arr[52,5] = 0
$arr[53,0] = "Chicago, IL"
$arr[53,1] = "US"
$arr[53,2] = "2/1/20 19:43"
$arr[53,3] = 2
$arr[53,4] = 0
$arr[53,5] = 0
$arr[54,0] = "San Benito, CA"
$arr[54,1] = "US"
$arr[54,2] = "2/3/20 3:53"
$arr[54,3] = 2
$arr[54,4] = 0
$arr[54,5] = 0
$arr[55,0] = "Santa Clara, CA"
$arr[55,1] = "US"
$arr[55,2] = "2/3/20 0:43"
$arr[55,3] = 2
$arr[55,4] = 0
$arr[55,5] = 0
$arr[56,0] = ""
$arr[56,1] = "Belgium"
$arr[56,2] = "2/4/20 15:43"
$arr[56,3] = 1
$arr[56,4] = 0
$arr[56,5] = 0
$arr[57,0] = ""
$arr[57,1] = "Cambodia"
$arr[57,2] = "1/31/20 8:15"
$arr[57,3] = 1
$arr[57,4] = 0
$arr[57,5] = 0
$arr[58,0] = "London, ON"
$arr[58,1] = "Canada"
$arr[58,2] = "2/4/20 0:03"
$arr[58,3] = 1
$arr[58,4] = 0
$arr[58,5] = 0
$arr[59,0] = ""
$arr[59,1] = "Finland"
$arr[59,2] = "1/31/20 8:15"
$arr[59,3] = 1
$arr[59,4] = 0
$arr[59,5] = 0
$arr[60,0] = "Tibet"
$arr[60,1] = "Mainland China"
$arr[60,2] = "2/1/20 1:52"
$arr[60,3] = 1
$arr[60,4] = 0
$arr[60,5] = 0
$arr[61,0] = ""
$arr[61,1] = "Nepal"
$arr[61,2] = "1/31/20 8:15"
$arr[61,3] = 1
$arr[61,4] = 0
$arr[61,5] = 0
$arr[62,0] = ""
$arr[62,1] = "Spain"
$arr[62,2] = "2/1/20 23:43"
$arr[62,3] = 1
$arr[62,4] = 0
$arr[62,5] = 0
$arr[63,0] = ""
$arr[63,1] = "Sri Lanka"
$arr[63,2] = "1/31/20 8:15"
$arr[63,3] = 1
$arr[63,4] = 0
$arr[63,5] = 0
$arr[64,0] = ""
$arr[64,1] = "Sweden"
$arr[64,2] = "2/1/20 2:13"
$arr[64,3] = 1
$arr[64,4] = 0
$arr[64,5] = 0
$arr[65,0] = "Boston, MA"
$arr[65,1] = "US"
$arr[65,2] = "2/1/20 19:43"
$arr[65,3] = 1
$arr[65,4] = 0
$arr[65,5] = 0
$arr[66,0] = "Los Angeles, CA"
$arr[66,1] = "US"
$arr[66,2] = "2/1/20 19:53"
$arr[66,3] = 1
$arr[66,4] = 0
$arr[66,5] = 0
$arr[67,0] = "Madison, WI"
$arr[67,1] = "US"
$arr[67,2] = "2/5/20 21:53"
$arr[67,3] = 1
$arr[67,4] = 0
$arr[67,5] = 0
$arr[68,0] = "Orange, CA"
$arr[68,1] = "US"
$arr[68,2] = "2/1/20 19:53"
$arr[68,3] = 1
$arr[68,4] = 0
$arr[68,5] = 0
$arr[69,0] = "Seattle, WA"
$arr[69,1] = "US"
$arr[69,2] = "2/1/20 19:43"
$arr[69,3] = 1
$arr[69,4] = 0
$arr[69,5] = 0
$arr[70,0] = "Tempe, AZ"
$arr[70,1] = "US"
$arr[70,2] = "2/1/20 19:43"
$arr[70,3] = 1
$arr[70,4] = 0
$arr[70,5] = 0

$ws.Range("A1977:F2047").Value = $arr

